# STA_metrics.xlsx update:
#  - insert two new columns (sum_SASA, max_SASA) after the SASA column
#  - replace the data block with a single updated row
#  - remove the now-obsolete rows 3-7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D so the old D (flexibility) moves to F, etc.
$ws.Range("D:E").EntireColumn.Insert()

# New header labels for the inserted columns
$ws.Range("D1").Value = "sum_SASA"
$ws.Range("E1").Value = "max_SASA"

# Delete the now-stale rows 3 through 7, leaving only the header + 1 data row
$ws.Range("3:7").EntireRow.Delete()

# Update row 2 with the new data values
$ws.Range("A2").Value = "GlcNAc(b1-4)GlcNAc(b1-4)GlcNAc"
$ws.Range("B2").Value = 1.922476871100382
$ws.Range("C2").Value = 3.321409580704815
$ws.Range("D2").Value = 3.321409580704815
$ws.Range("E2").Value = 3.321409580704815
$ws.Range("F2").Value = 1.234
$ws.Range("G2").Value = 0.371
$ws.Range("H2").Value = 8.91
$ws.Range("I2").Value = "['4C1']"
$ws.Range("J2").Value = "['GlcNAc(b1-4)']"
$ws.Range("K2").Value = "['GlcNAc(b1-4)']"
$ws.Range("L2").Value = "N"
